$wb = $excel.ActiveWorkbook

# Sheet "OFF" - row 3 (R) updates
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 216
$wsOff.Range("C3").Value = 156
$wsOff.Range("D3").Value = 51
$wsOff.Range("E3").Value = 29

# Sheet "DEF" - row 3 (R) updates
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 184
$wsDef.Range("C3").Value = 124
$wsDef.Range("D3").Value = 57
$wsDef.Range("E3").Value = 24
$wsDef.Range("F3").Value = 4
